$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = Get-Date -Year 2023 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
$ws.Range("C6").Value = $newDate
